$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Date" column (BF) holds a mis-derived date string ("5-27-2012-13")
# for every data row (rows 2-31). Replace it with the correct ISO date
# ("2013-05-27"). Writing that literal text straight into .Value would be
# auto-recognized as a real date (like typing it into Excel) and silently
# converted to a date serial number, which is not what the source data
# looks like. So: force the cell(s) to Text first, write the literal
# string, then drop the format back to the workbook's default "Normal"
# style so no stray per-cell formatting is introduced.
$range = $ws.Range("BF2:BF31")
$range.NumberFormat = "@"

for ($r = 2; $r -le 31; $r++) {
    $ws.Range("BF$r").Value = "2013-05-27"
}

$range.Style = "Normal"
